# Auto-generated Excel COM-interop script applying the Durandal_Profits value updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
# Row 19
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("N19").ClearContents()

# Row 137
$ws.Range("H137").Value = 886.45
$ws.Range("I137").Value = 774.7917
$ws.Range("J137").Value = 1053.9375
$ws.Range("K137").Value = 2324.3751
$ws.Range("L137").Value = 3161.8125
$ws.Range("M137").Value = 225.6248999999998
$ws.Range("N137").Value = -8261.8125

# Row 138
$ws.Range("H138").Value = 4786.4443
$ws.Range("I138").Value = 3431.8333
$ws.Range("J138").Value = 5689.5186
$ws.Range("K138").Value = 10295.4999
$ws.Range("L138").Value = 17068.5558
$ws.Range("M138").Value = -5155.499899999999
$ws.Range("N138").Value = -27348.5558

$ws = $wb.Worksheets("ARM")
# Row 2
$ws.Range("H2").Value = 2163.76
$ws.Range("I2").Value = 1229.05
$ws.Range("J2").Value = 5902.6
$ws.Range("K2").Value = 1229.05
$ws.Range("L2").Value = 5902.6
$ws.Range("M2").Value = -1116.05
$ws.Range("N2").Value = -6128.6

# Row 61
$ws.Range("H61").Value = 1602.4445
$ws.Range("I61").Value = 1365.25
$ws.Range("K61").Value = 1365.25
$ws.Range("M61").Value = -1153.25

# Row 116
$ws.Range("H116").Value = 2163.76
$ws.Range("I116").Value = 1229.05
$ws.Range("J116").Value = 5902.6
$ws.Range("K116").Value = 1229.05
$ws.Range("L116").Value = 5902.6
$ws.Range("M116").Value = 1064.95
$ws.Range("N116").Value = -10490.6

# Row 136
$ws.Range("H136").Value = 1602.4445
$ws.Range("I136").Value = 1365.25
$ws.Range("K136").Value = 4095.75
$ws.Range("M136").Value = -1545.75

$ws = $wb.Worksheets("BSM")
# Row 3
$ws.Range("H3").Value = 2163.76
$ws.Range("I3").Value = 1229.05
$ws.Range("J3").Value = 5902.6
$ws.Range("K3").Value = 1229.05
$ws.Range("L3").Value = 5902.6
$ws.Range("M3").Value = -1115.05
$ws.Range("N3").Value = -6130.6

# Row 20
$ws.Range("H20").Value = 2575.4243
$ws.Range("I20").Value = 1984.8889
$ws.Range("J20").Value = 3284.0667
$ws.Range("K20").Value = 1984.8889
$ws.Range("L20").Value = 3284.0667
$ws.Range("M20").Value = -1737.8889
$ws.Range("N20").Value = -3778.0667

# Row 94
$ws.Range("H94").Value = 1116.88
$ws.Range("I94").Value = 1155.6666
$ws.Range("J94").Value = 1017.1429
$ws.Range("K94").Value = 1155.6666
$ws.Range("L94").Value = 1017.1429
$ws.Range("M94").Value = -704.6666
$ws.Range("N94").Value = -1919.1429

$ws = $wb.Worksheets("CRP")
# Row 31
$ws.Range("H31").Value = 3708.681
$ws.Range("I31").Value = 1898.2162
$ws.Range("J31").Value = 10407.4
$ws.Range("K31").Value = 1898.2162
$ws.Range("L31").Value = 10407.4
$ws.Range("M31").Value = -1603.2162
$ws.Range("N31").Value = -10997.4

# Row 34
$ws.Range("H34").Value = 3708.681
$ws.Range("I34").Value = 1898.2162
$ws.Range("J34").Value = 10407.4
$ws.Range("K34").Value = 1898.2162
$ws.Range("L34").Value = 10407.4
$ws.Range("M34").Value = -1696.2162
$ws.Range("N34").Value = -10811.4

# Row 58
$ws.Range("H58").Value = 1364.5151
$ws.Range("I58").Value = 1141.6666
$ws.Range("J58").Value = 1754.5
$ws.Range("K58").Value = 1141.6666
$ws.Range("L58").Value = 1754.5
$ws.Range("M58").Value = -938.6666
$ws.Range("N58").Value = -2160.5

# Row 99
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").ClearContents()

# Row 126
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()

# Row 132
$ws.Range("H132").Value = 1264.8422
$ws.Range("I132").Value = 979.1163
$ws.Range("J132").Value = 2142.4285
$ws.Range("K132").Value = 2937.3489
$ws.Range("L132").Value = 6427.2855
$ws.Range("M132").Value = -407.3489
$ws.Range("N132").Value = -11487.2855

# Row 136
$ws.Range("H136").Value = 1364.5151
$ws.Range("I136").Value = 1141.6666
$ws.Range("J136").Value = 1754.5
$ws.Range("K136").Value = 3424.9998
$ws.Range("L136").Value = 5263.5
$ws.Range("M136").Value = -874.9998000000001
$ws.Range("N136").Value = -10363.5

$ws = $wb.Worksheets("CUL")
# Row 5
$ws.Range("H5").Value = 808.775
$ws.Range("I5").Value = 451.13043
$ws.Range("J5").Value = 1292.6471
$ws.Range("K5").Value = 1353.39129
$ws.Range("L5").Value = 3877.9413
$ws.Range("M5").Value = -1241.39129
$ws.Range("N5").Value = -4101.9413

# Row 12
$ws.Range("H12").Value = 92.92308
$ws.Range("I12").Value = 15.571428
$ws.Range("J12").Value = 121.42105
$ws.Range("K12").Value = 46.714284
$ws.Range("L12").Value = 364.26315
$ws.Range("M12").Value = 126.285716
$ws.Range("N12").Value = -710.26315

# Row 38
$ws.Range("H38").Value = 1702.8889
$ws.Range("I38").Value = 1035.2
$ws.Range("J38").Value = 2537.5
$ws.Range("K38").Value = 3105.6
$ws.Range("L38").Value = 7612.5
$ws.Range("M38").Value = -2758.6
$ws.Range("N38").Value = -8306.5

# Row 131
$ws.Range("H131").Value = 6579778
$ws.Range("J131").Value = 7247243
$ws.Range("L131").Value = 21741729
$ws.Range("N131").Value = -21751809

# Row 135
$ws.Range("H135").Value = 808.775
$ws.Range("I135").Value = 451.13043
$ws.Range("J135").Value = 1292.6471
$ws.Range("K135").Value = 4060.17387
$ws.Range("L135").Value = 11633.8239
$ws.Range("M135").Value = -1525.17387
$ws.Range("N135").Value = -16703.8239

$ws = $wb.Worksheets("GSM")
# Row 2
$ws.Range("H2").Value = 35.545456
$ws.Range("I2").Value = 21.285715
$ws.Range("J2").Value = 60.5
$ws.Range("K2").Value = 21.285715
$ws.Range("L2").Value = 60.5
$ws.Range("M2").Value = 91.714285
$ws.Range("N2").Value = -286.5

# Row 70
$ws.Range("H70").Value = 12237206
$ws.Range("I70").Value = 29617204
$ws.Range("J70").Value = 6837.593
$ws.Range("K70").Value = 29617204
$ws.Range("L70").Value = 6837.593
$ws.Range("M70").Value = -29616934
$ws.Range("N70").Value = -7377.593

# Row 73
$ws.Range("H73").Value = 12237206
$ws.Range("I73").Value = 29617204
$ws.Range("J73").Value = 6837.593
$ws.Range("K73").Value = 29617204
$ws.Range("L73").Value = 6837.593
$ws.Range("M73").Value = -29616268
$ws.Range("N73").Value = -8709.593000000001

# Row 97
$ws.Range("H97").Value = 941.3871
$ws.Range("I97").Value = 594.7368
$ws.Range("J97").Value = 1490.25
$ws.Range("K97").Value = 594.7368
$ws.Range("L97").Value = 1490.25
$ws.Range("M97").Value = -98.73680000000002
$ws.Range("N97").Value = -2482.25

# Row 102
$ws.Range("H102").Value = 3014
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 3014
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 3014
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = -6258

# Row 134
$ws.Range("H134").Value = 15362.875
$ws.Range("J134").Value = 15362.875
$ws.Range("L134").Value = 46088.625
$ws.Range("N134").Value = -51158.625

# Row 136
$ws.Range("H136").Value = 12428.533
$ws.Range("J136").Value = 12428.533
$ws.Range("L136").Value = 37285.599
$ws.Range("N136").Value = -42385.599

$ws = $wb.Worksheets("WVR")
# Row 41
$ws.Range("H41").Value = 6780.5
$ws.Range("J41").Value = 7002.1816
$ws.Range("L41").Value = 7002.1816
$ws.Range("N41").Value = -7782.1816

# Row 132
$ws.Range("H132").Value = 20834542
$ws.Range("I132").Value = 24510876
$ws.Range("J132").Value = 1980.8889
$ws.Range("K132").Value = 73532628
$ws.Range("L132").Value = 5942.6667
$ws.Range("M132").Value = -73530098
$ws.Range("N132").Value = -11002.6667

# Row 136
$ws.Range("H136").Value = 928
$ws.Range("I136").Value = 808.5405
$ws.Range("J136").Value = 1370
$ws.Range("K136").Value = 2425.6215
$ws.Range("L136").Value = 4110
$ws.Range("M136").Value = 124.3785000000003
$ws.Range("N136").Value = -9210
